# Update results in pl_mw.xlsx (active sheet) for rows 2-25, columns C,D,E,F,G,K,L,M,N
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters in the order the data columns appear below
$cols = @("C","D","E","F","G","K","L","M","N")

# New values per row (rows 2..25), in the same column order as $cols
$newValues = @(
    @(0.3259446694196271,0.2238569003226374,0.1048970261287216,3.474901115801543,0.002615568426622285,3.546465783198983,0.06917241373459504,0.9129259863411221,2.662435476810145),
    @(0.32325184456964,0.2252201985199136,0.1047270739702633,3.434296381402092,0.002622198047681925,3.392826100289255,0.06941196327400867,0.8845560840256681,2.679212313847131),
    @(0.3217662373620414,0.2261338222959139,0.1046681294139464,3.411345103296512,0.002626478611197891,3.30082909010315,0.06958029540156829,0.867718929957654,2.690291627598555),
    @(0.3212029201395836,0.226525324349538,0.1046555100902928,3.402487323719185,0.002628275967848862,3.263923237860865,0.06965423542116866,0.8610031790066088,2.695001513394743),
    @(0.3211119199862225,0.2265914905284845,0.1046541027804118,3.401046312678318,0.00262857762416058,3.257830177199992,0.06966683584386679,0.8598967988428612,2.695795341280245),
    @(0.3217584700365705,0.2261390245918236,0.1046679130842438,3.41122364349819,0.002626502636143788,3.300329007228584,0.06958127094707223,0.8676277708586184,2.690354358263249),
    @(0.3249812854281089,0.2243110423242172,0.1048289918242347,3.460487846171333,0.002617810864966889,3.493002804488867,0.06925060057131027,0.9030227094187211,2.668058162260323),
    @(0.3326386447172069,0.2213364108615892,0.1055060711329112,3.572947290491157,0.00260242309108413,3.889611802509194,0.06877076642233249,0.9770941129468298,2.630543259273338),
    @(0.3390894383793182,0.2195265515700768,0.1062252354952697,3.665445243828742,0.002592114938967646,4.192794085137962,0.06852109911072901,1.034424629316959,2.606810656672451),
    @(0.3422054132770427,0.2187855907370064,0.1066009014630893,3.70971609779059,0.002587639300028006,4.333359409518664,0.0684298808876953,1.061153244894584,2.596855593491981),
    @(0.3434116142685184,0.2185169223020793,0.1067501596563147,3.72679934666192,0.002585974995366796,4.386974229592113,0.06839855587813304,1.071369112308389,2.593207668693296),
    @(0.3431506673904039,0.2185742536137099,0.1067177024198713,3.723105926864434,0.002586332078777343,4.375410075754132,0.06840515914413814,1.069164730644331,2.593987881044882),
    @(0.3423041211985662,0.2187632479447075,0.1066130405523751,3.711115136823281,0.002587501765972548,4.337762570746804,0.06842723926418159,1.061991814058743,2.596553029401903),
    @(0.3417890102156207,0.2188805666816478,0.106549844726139,3.703812058085362,0.002588222202470766,4.314752785672226,0.06844118304887914,1.057610510284817,2.598140151315576),
    @(0.3388894668163687,0.2195766386340026,0.1062016632178384,3.662596465175994,0.002592411713832305,4.183661505559769,0.06852751051566841,1.032690992888945,2.607478243812011),
    @(0.3371572785499382,0.2200248016712862,0.1060005085611557,3.637876183441477,0.002595036413082467,4.103922705935929,0.06858619737931448,1.017570550268061,2.613422904055724),
    @(0.33617804166974,0.2202903200726709,0.1058893745617837,3.623863941399833,0.002596566186440644,4.058307805817549,0.06862205666136489,1.008934664082531,2.616921220858671),
    @(0.3358494164418744,0.2203815483166522,0.1058525297594883,3.619154948983606,0.00259708760182073,4.04290594894195,0.06863455931710405,1.006021145622213,2.618119249936896),
    @(0.3373399049695536,0.2199762916424461,0.1060214491942766,3.640486332708946,0.002594754928810809,4.112385255059223,0.0685797322737649,1.019173827248295,2.612781891213473),
    @(0.3425520587056212,0.2187074116871628,0.106643592050478,3.714628441310367,0.002587157373544966,4.348810042909633,0.06842066645827316,1.064096105253,2.595796269113464),
    @(0.346111575714815,0.2179476249371106,0.1070910217931633,3.764944686653848,0.002582369753834405,4.505578483641671,0.06833546231561272,1.094005688007684,2.585405800326441),
    @(0.3441977365143032,0.2183467528997269,0.1068484761411881,3.737918594051422,0.002584908788978986,4.421700459241038,0.06837922042466715,1.077991678845379,2.590886064488643),
    @(0.3372572878057554,0.2199981985384483,0.1060119678805052,3.639305662225581,0.002594882123028561,4.108558626833997,0.06858264854683682,1.018448808081885,2.613071441973602),
    @(0.3304229692862322,0.2220754420058384,0.1052840807087136,3.540805448120324,0.002606409844408522,3.780272996659505,0.06888251644074828,0.9565497805366618,2.640022830139557)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellRef = "{0}{1}" -f $cols[$j], $row
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}
